$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A (the numeric index column) - this shifts old B,C,D,E left to A,B,C,D
$ws.Range("A1:A4").EntireColumn.Delete()

# Resize the remaining columns (now A,B,C,D) to their new widths
$ws.Columns.Item(1).ColumnWidth = 26
$ws.Columns.Item(2).ColumnWidth = 32.333333333333336
$ws.Columns.Item(3).ColumnWidth = 38.166666666666664
$ws.Columns.Item(4).ColumnWidth = 59.166666666666664

# Add wrapText to the header row style (bold/bordered header cells)
$ws.Range("A1:D1").WrapText = $true

# Update row heights for the data rows to match the new wrapped layout
$ws.Rows.Item(2).RowHeight = 129.6
$ws.Rows.Item(3).RowHeight = 374.4
$ws.Rows.Item(4).RowHeight = 409.6

# Update the selected cell
$ws.Range("A2").Select()
